$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking "Price" (D) values are plain numeric-looking strings (e.g. "1.00",
# "219.01") that must stay text, matching the inlineStr cells in the source sheet.
# Forcing NumberFormat to Text before the write keeps Excel from auto-converting
# them to floats, and resetting the style back to Normal afterwards keeps the
# cell style index identical to the original (no explicit "s" attribute).
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "84.049.61"
$ws.Range("E2").Value = "  +5.51%  "
# Row 3
Set-TextValue "D3" "3.246.10"
$ws.Range("E3").Value = "  +1.26%  "
# Row 4
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  +0.15%  "
# Row 5
Set-TextValue "D5" "219.01"
$ws.Range("E5").Value = "  +3.65%  "
# Row 6
Set-TextValue "D6" "625.74"
$ws.Range("E6").Value = "  -1.60%  "
# Row 7
Set-TextValue "D7" "0.309"
$ws.Range("E7").Value = "  +25.67%  "
# Row 8
Set-TextValue "D8" "1.00"
$ws.Range("E8").Value = "  +0.05%  "
# Row 9
Set-TextValue "D9" "0.593"
$ws.Range("E9").Value = "  -1.34%  "
# Row 10
Set-TextValue "D10" "3.243.99"
$ws.Range("E10").Value = "  +1.38%  "
# Row 11
Set-TextValue "D11" "0.612"
$ws.Range("E11").Value = "  +2.21%  "
# Row 12
Set-TextValue "D12" "0.0000280"
$ws.Range("E12").Value = "  +7.28%  "
# Row 13
$ws.Range("E13").Value = "  -0.06%  "
# Row 14
Set-TextValue "D14" "5.40"
$ws.Range("E14").Value = "  -0.24%  "
# Row 15
Set-TextValue "D15" "3.835.25"
$ws.Range("E15").Value = "  +1.28%  "
# Row 16
Set-TextValue "D16" "32.69"
$ws.Range("E16").Value = "  +1.28%  "
# Row 17
Set-TextValue "D17" "83.424.72"
$ws.Range("E17").Value = "  +4.90%  "
# Row 18
Set-TextValue "D18" "3.232.31"
$ws.Range("E18").Value = "  +1.06%  "
# Row 19
Set-TextValue "D19" "3.27"
$ws.Range("E19").Value = "  +8.48%  "
# Row 20
Set-TextValue "D20" "14.42"
$ws.Range("E20").Value = "  -0.80%  "
# Row 21
Set-TextValue "D21" "447.45"
$ws.Range("E21").Value = "  +1.60%  "
# Row 22
Set-TextValue "D22" "9.12"
$ws.Range("E22").Value = "  -2.54%  "
# Row 23
Set-TextValue "D23" "5.24"
$ws.Range("E23").Value = "  +0.26%  "
# Row 24
Set-TextValue "D24" "7.42"
$ws.Range("E24").Value = "  +6.31%  "
# Row 25
Set-TextValue "D25" "5.23"
$ws.Range("E25").Value = "  +8.37%  "
# Row 26
Set-TextValue "D26" "11.97"
$ws.Range("E26").Value = "  +9.94%  "
# Row 27
$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D27" "78.60"
$ws.Range("E27").Value = "  +1.70%  "
# Row 28
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue "D28" "3.394.47"
$ws.Range("E28").Value = "  +0.70%  "
# Row 29
Set-TextValue "D29" "1.01"
$ws.Range("E29").Value = "  +1.00%  "
# Row 30
Set-TextValue "D30" "0.0000125"
$ws.Range("E30").Value = "  -0.15%  "
# Row 31
Set-TextValue "D31" "9.20"
$ws.Range("E31").Value = "  +0.33%  "
# Row 32
$ws.Range("E32").Value = "  -0.11%  "
# Row 33
Set-TextValue "D33" "570.89"
$ws.Range("E33").Value = "  +2.24%  "
# Row 34
$ws.Range("B34").Value = "Cronos"
$ws.Range("C34").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D34" "0.150"
$ws.Range("E34").Value = "  +23.33%  "
# Row 35
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D35" "1.49"
$ws.Range("E35").Value = "  -2.81%  "
# Row 36
$ws.Range("E36").Value = "  -2.15%  "
# Row 37
$ws.Range("E37").Value = "  -2.12%  "
# Row 38
Set-TextValue "D38" "23.24"
$ws.Range("E38").Value = "  +0.60%  "
# Row 39
Set-TextValue "D39" "6.22"
$ws.Range("E39").Value = "  +10.09%  "
# Row 40
Set-TextValue "D40" "1.00"
$ws.Range("E40").Value = "  +0.28%  "
# Row 41
$ws.Range("E41").Value = "  -1.20%  "
# Row 42
$ws.Range("B42").Value = "WhiteBITCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D42" "20.93"
$ws.Range("E42").Value = "  +3.82%  "
# Row 43
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D43" "2.04"
$ws.Range("E43").Value = "  +11.54%  "
# Row 44
$ws.Range("E44").Value = "  +12.37%  "
# Row 45
Set-TextValue "D45" "160.19"
$ws.Range("E45").Value = "  -1.96%  "
# Row 46
$ws.Range("E46").Value = "  +0.08%  "
# Row 47
Set-TextValue "D47" "188.77"
$ws.Range("E47").Value = "  -1.95%  "
# Row 48
Set-TextValue "D48" "45.03"
$ws.Range("E48").Value = "  +3.92%  "
# Row 49
Set-TextValue "D49" "1.32"
$ws.Range("E49").Value = "  -1.70%  "
# Row 50
Set-TextValue "D50" "0.784"
$ws.Range("E50").Value = "  -2.00%  "
# Row 51
Set-TextValue "D51" "26.04"
$ws.Range("E51").Value = "  +1.63%  "
